$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 486.5
$ws.Range("I11").Value = 486.5
$ws.Range("K11").Value = 486.5
$ws.Range("M11").Value = -346.5

$ws.Range("H107").Value = 1127.1666
$ws.Range("J107").Value = 1622.5
$ws.Range("L107").Value = 1622.5
$ws.Range("N107").Value = -5462.5

$ws.Range("H127").Value = 3200
$ws.Range("I127").Value = 3200
$ws.Range("K127").Value = 9600
$ws.Range("M127").Value = -4640

$ws.Range("H138").Value = 2912.0615
$ws.Range("J138").Value = 2956.0164
$ws.Range("L138").Value = 8868.049199999999
$ws.Range("N138").Value = -19148.0492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18034.535
$ws.Range("I32").Value = 6033.34
$ws.Range("K32").Value = 6033.34
$ws.Range("M32").Value = -5746.34

$ws.Range("H122").Value = 1904.7273
$ws.Range("I122").Value = 1904.7273
$ws.Range("K122").Value = 5714.1819
$ws.Range("M122").Value = -3264.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1129.7222
$ws.Range("I20").Value = 1061.2142
$ws.Range("J20").Value = 1369.5
$ws.Range("K20").Value = 1061.2142
$ws.Range("L20").Value = 1369.5
$ws.Range("M20").Value = -814.2141999999999
$ws.Range("N20").Value = -1863.5

$ws.Range("H80").Value = 1704.5
$ws.Range("I80").Value = 353.25
$ws.Range("J80").Value = 2245
$ws.Range("K80").Value = 353.25
$ws.Range("L80").Value = 2245
$ws.Range("M80").Value = 644.75
$ws.Range("N80").Value = -4241

$ws.Range("H83").Value = 1704.5
$ws.Range("I83").Value = 353.25
$ws.Range("J83").Value = 2245
$ws.Range("K83").Value = 1766.25
$ws.Range("L83").Value = 11225
$ws.Range("M83").Value = 3225.75
$ws.Range("N83").Value = -21209

$ws.Range("H99").Value = 1624.5
$ws.Range("I99").Value = 1499.3334
$ws.Range("K99").Value = 1499.3334
$ws.Range("M99").Value = -1.333399999999983

$ws.Range("H104").Value = 99000
$ws.Range("J104").Value = 99000
$ws.Range("L104").Value = 99000
$ws.Range("N104").Value = -105988

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1261.4
$ws.Range("I8").Value = 1469.6666
$ws.Range("J8").Value = 949
$ws.Range("K8").Value = 1469.6666
$ws.Range("L8").Value = 949
$ws.Range("M8").Value = -1329.6666
$ws.Range("N8").Value = -1229

$ws.Range("H58").Value = 2653.7058
$ws.Range("I58").Value = 2569.5625
$ws.Range("K58").Value = 2569.5625
$ws.Range("M58").Value = -2366.5625

$ws.Range("H99").Value = 4980.6665
$ws.Range("I99").Value = 4980.6665
$ws.Range("K99").Value = 4980.6665
$ws.Range("M99").Value = -3482.6665

$ws.Range("H122").Value = 2888.0588
$ws.Range("J122").Value = 4999
$ws.Range("L122").Value = 14997
$ws.Range("N122").Value = -19897

$ws.Range("H126").Value = 4980.6665
$ws.Range("I126").Value = 4980.6665
$ws.Range("K126").Value = 14941.9995
$ws.Range("M126").Value = -12471.9995

$ws.Range("H132").Value = 2202
$ws.Range("I132").Value = 1842.5883
$ws.Range("K132").Value = 5527.7649
$ws.Range("M132").Value = -2997.7649

$ws.Range("H136").Value = 2653.7058
$ws.Range("I136").Value = 2569.5625
$ws.Range("K136").Value = 7708.6875
$ws.Range("M136").Value = -5158.6875

$ws.Range("H141").Value = 210834.22
$ws.Range("J141").Value = 210834.22
$ws.Range("L141").Value = 210834.22
$ws.Range("N141").Value = -221194.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 669.85
$ws.Range("I107").Value = 337.22223
$ws.Range("K107").Value = 1011.66669
$ws.Range("M107").Value = 908.33331

$ws.Range("H108").Value = 1863.5
$ws.Range("I108").Value = 727
$ws.Range("K108").Value = 2181
$ws.Range("M108").Value = 699

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1210.4814
$ws.Range("I97").Value = 1218.1538
$ws.Range("K97").Value = 1218.1538
$ws.Range("M97").Value = -722.1538

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents() | Out-Null

$ws.Range("H122").Value = 2946.2083
$ws.Range("I122").Value = 1779.8
$ws.Range("K122").Value = 5339.4
$ws.Range("M122").Value = -2889.4

$ws.Range("H132").Value = 3367.4375
$ws.Range("I132").Value = 3995.0527
$ws.Range("J132").Value = 2450.1538
$ws.Range("K132").Value = 11985.1581
$ws.Range("L132").Value = 7350.4614
$ws.Range("M132").Value = -9455.158100000001
$ws.Range("N132").Value = -12410.4614

$ws.Range("H134").Value = 109499.5
$ws.Range("J134").Value = 109499.5
$ws.Range("L134").Value = 328498.5
$ws.Range("N134").Value = -333568.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 568.1818
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 605.55554
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 605.55554
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1195.55554

$ws.Range("H27").Value = 568.1818
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 605.55554
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 605.55554
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -819.55554

$ws.Range("H46").Value = 1873.5807
$ws.Range("I46").Value = 1428.7333
$ws.Range("K46").Value = 1428.7333
$ws.Range("M46").Value = -1240.7333

$ws.Range("H61").Value = 3184.8948
$ws.Range("I61").Value = 3200.7778
$ws.Range("J61").Value = 2899
$ws.Range("K61").Value = 3200.7778
$ws.Range("L61").Value = 2899
$ws.Range("M61").Value = -2998.7778
$ws.Range("N61").Value = -3303

$ws.Range("H109").Value = 90285
$ws.Range("J109").Value = 90285
$ws.Range("L109").Value = 90285
$ws.Range("N109").Value = -93059

$ws.Range("H113").Value = 3184.8948
$ws.Range("I113").Value = 3200.7778
$ws.Range("J113").Value = 2899
$ws.Range("K113").Value = 3200.7778
$ws.Range("L113").Value = 2899
$ws.Range("M113").Value = -1030.7778
$ws.Range("N113").Value = -7239

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents() | Out-Null

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 95000
$ws.Range("J16").Value = 95000
$ws.Range("L16").Value = 95000
$ws.Range("N16").Value = -95584

$ws.Range("H41").Value = 9079
$ws.Range("J41").Value = 6662
$ws.Range("L41").Value = 6662
$ws.Range("N41").Value = -7442

$ws.Range("H122").Value = 22861.588
$ws.Range("I122").Value = 23178.066
$ws.Range("J122").Value = 20488
$ws.Range("K122").Value = 69534.198
$ws.Range("L122").Value = 61464
$ws.Range("M122").Value = -67084.198
$ws.Range("N122").Value = -66364

$ws.Range("H132").Value = 2301.975
$ws.Range("I132").Value = 2318.5715
$ws.Range("J132").Value = 2263.25
$ws.Range("K132").Value = 6955.7145
$ws.Range("L132").Value = 6789.75
$ws.Range("M132").Value = -4425.7145
$ws.Range("N132").Value = -11849.75

$ws.Range("H135").Value = 88750
$ws.Range("J135").Value = 88750
$ws.Range("L135").Value = 88750
$ws.Range("N135").Value = -98890

$ws.Range("H136").Value = 1735.1428
$ws.Range("I136").Value = 1750.1111
$ws.Range("K136").Value = 5250.3333
$ws.Range("M136").Value = -2700.3333

$ws.Range("H139").Value = 110000
$ws.Range("J139").Value = 110000
$ws.Range("L139").Value = 110000
$ws.Range("N139").Value = -120280
